# Insert a new weekly record at row 842, pushing the rest of the
# "Terminal La Palmera de La Serena - Zapallo / Camote" table down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row at row 842 - this shifts rows 842:954 down to 843:955
$ws.Rows.Item(842).Insert()

# Fill the new row 842 with the new weekly data.
$ws.Cells.Item(842, 1).Value = 8
$ws.Cells.Item(842, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(842, 3).Value = "Coquimbo"
$ws.Cells.Item(842, 4).Value = 45127
$ws.Cells.Item(842, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(842, 5).Value = 4
$ws.Cells.Item(842, 6).Value = 100112045
$ws.Cells.Item(842, 7).Value = "Zapallo"
$ws.Cells.Item(842, 8).Value = "Camote"
$ws.Cells.Item(842, 9).Value = "1a (guarda)"
$ws.Cells.Item(842, 10).Value = 1000
$ws.Cells.Item(842, 11).Value = 600
$ws.Cells.Item(842, 12).Value = 700
$ws.Cells.Item(842, 13).Value = 650
$ws.Cells.Item(842, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(842, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(842, 16).Value = 650
$ws.Cells.Item(842, 17).Value = 1
$ws.Cells.Item(842, 18).Value = "Hortaliza"
